# Regenerate save_data to use K (strikeouts, column G) instead of the old
# Strike# counts. New K values were recalculated (calc and write s_vals)
# and the std/mean derived from them changed as a result.
#
# Column G ("K") holds per-row counts; rows 2-45 get new values below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 3
    4  = 3
    5  = 2
    6  = 0
    7  = 4
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 2
    25 = 0
    26 = 1
    27 = 2
    28 = 0
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 2
    36 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 0
    44 = 1
    45 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
